# Renumber the "Session" column (D) for every "Surgery Seminar/Slide" row (column C)
# from the sheet-wide numbering (16-22) down to a per-subject numbering (1-7),
# storing the value as a literal number instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $subject = $ws.Cells.Item($r, 3).Value()
    if ($subject -eq "Surgery Seminar/Slide") {
        $oldSession = $ws.Cells.Item($r, 4).Value()
        $newSession = [int]$oldSession - 15
        $ws.Cells.Item($r, 4).Value = $newSession
    }
}

$wb.Save()
